$d = $word.ActiveDocument

function Merge-And-Unbold([string]$fullText, [int]$firstRunLen) {
    # Locate the full placeholder text (Word's Find matches across run boundaries).
    $whole = $d.Content
    $whole.Find.ClearFormatting()
    $whole.Find.Text = $fullText
    $whole.Find.Execute($fullText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($whole.Find.Found) {
        $start = $whole.Start
        $end = $whole.End
        if ($firstRunLen -lt ($end - $start)) {
            # Everything after the first run's original text is in trailing run(s);
            # delete that trailing range so only one run (the first) remains, then
            # append the removed text back onto the (now sole) run.
            $tailStart = $start + $firstRunLen
            $tail = $d.Range($tailStart, $end)
            $tailText = $tail.Text
            $tail.Delete()

            $remain = $d.Range($start, $start + $firstRunLen)
            $remain.InsertAfter($tailText)
        }
    }

    # Re-find fresh (post edit) so the Font getter/setter targets the live run.
    $final = $d.Content
    $final.Find.ClearFormatting()
    $final.Find.Text = $fullText
    $final.Find.Execute($fullText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($final.Find.Found) {
        $final.Font.Bold = 0
        $final.Font.BoldBi = 0
    }
}

# {{Certification_1}} is already a single run -> just drop bold.
Merge-And-Unbold "{{Certification_1}}" 20

# {{Certification_2}} is split "{{Certification_" / "2" / "}}" -> merge + drop bold.
Merge-And-Unbold "{{Certification_2}}" 16

# {{Certification_3}} is split "{{Certification_" / "3" / "}}" -> merge + drop bold.
Merge-And-Unbold "{{Certification_3}}" 16

Write-Host "done"
